$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching style/format of existing headers (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill H2:H5 with 0 values (numeric), matching style/format of the corresponding G column cells
$ws.Range("H2:H5").Value = 0
$ws.Range("G2:G5").Copy()
$ws.Range("H2:H5").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
